# Insert two new rows at 684-685 (pushing the existing rows 684.. down to 686..)
# then populate them with the new weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("684:685").Insert()

# Row 684 (new) -------------------------------------------------------
$ws.Cells.Item(684, 1).Value = 11
$ws.Cells.Item(684, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(684, 3).Value = "Bíobío"
$ws.Cells.Item(684, 4).Value = 45077
$ws.Cells.Item(684, 5).Value = 8
$ws.Cells.Item(684, 6).Value = "Fruta"
$ws.Cells.Item(684, 7).Value = 100102
$ws.Cells.Item(684, 8).Value = "Cítricos"
$ws.Cells.Item(684, 9).Value = 100102003
$ws.Cells.Item(684, 10).Value = "Limón"
$ws.Cells.Item(684, 11).Value = "Sin especificar"
$ws.Cells.Item(684, 12).Value = "1a amarillo"
$ws.Cells.Item(684, 13).Value = 330
$ws.Cells.Item(684, 14).Value = 10000
$ws.Cells.Item(684, 15).Value = 11000
$ws.Cells.Item(684, 16).Value = 10545
$ws.Cells.Item(684, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(684, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(684, 19).Value = 659
$ws.Cells.Item(684, 20).Value = 16

# Row 685 (new) -------------------------------------------------------
$ws.Cells.Item(685, 1).Value = 11
$ws.Cells.Item(685, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(685, 3).Value = "Bíobío"
$ws.Cells.Item(685, 4).Value = 45077
$ws.Cells.Item(685, 5).Value = 8
$ws.Cells.Item(685, 6).Value = "Fruta"
$ws.Cells.Item(685, 7).Value = 100102
$ws.Cells.Item(685, 8).Value = "Cítricos"
$ws.Cells.Item(685, 9).Value = 100102003
$ws.Cells.Item(685, 10).Value = "Limón"
$ws.Cells.Item(685, 11).Value = "Sin especificar"
$ws.Cells.Item(685, 12).Value = "2a amarillo"
$ws.Cells.Item(685, 13).Value = 220
$ws.Cells.Item(685, 14).Value = 7000
$ws.Cells.Item(685, 15).Value = 8000
$ws.Cells.Item(685, 16).Value = 7545
$ws.Cells.Item(685, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(685, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(685, 19).Value = 472
$ws.Cells.Item(685, 20).Value = 16
